# Execute Tests on Android mobile browser
# Adds a "Locator Type" column (D) to the CapellaForm locator dictionary
# sheet, classifying each existing locator as "Xpath" or "CSS".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. New column D header + values
# ---------------------------------------------------------------------
$ws.Range("D1").Value = "Locator Type"

$locatorTypes = @(
    "Xpath",  # row 2  - CapellaPage_Button_AcceptCookies          (//div[...])
    "CSS",    # row 3  - CapellaPage_TextBox_FullName               (input#name)
    "CSS",    # row 4  - CapellaPage_ErrorMessage_...FullName       (#name + span.error-msg)
    "CSS",    # row 5  - CapellaPage_TextBox_Email                  (input#email)
    "CSS",    # row 6  - CapellaPage_ErrorMessage_...Email          (#email + span.error-msg)
    "CSS",    # row 7  - CapellaPage_TextBox_Password                (input#password)
    "CSS",    # row 8  - CapellaPage_TextBox_OrganizationName        (input#orgName)
    "CSS",    # row 9  - CapellaPage_ErrorMessage_...OrgName         (#orgName + span.error-msg)
    "Xpath",  # row 10 - CapellaPage_CheckBox_TermsOfService         (//input[...])
    "Xpath",  # row 11 - CapellaPage_CheckBox_MarketingOptIn         (//input[...])
    "CSS",    # row 12 - CapellaPage_Button_SignUp                   (button#formSubmit)
    "CSS",    # row 13 - CapellaPage_BrandLogo_ImageLink              (a[href='/'] img...)
    "CSS"     # row 14 - CapellaPage_InboxImage_SuccessfullConfirmation (img[src*=...])
)

for ($i = 0; $i -lt $locatorTypes.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 4).Value = $locatorTypes[$i]
}

# ---------------------------------------------------------------------
# 2. Header cell (D1) formatting: bold, 13pt, Helvetica Neue
# ---------------------------------------------------------------------
$ws.Range("D1").Font.Bold = $true
$ws.Range("D1").Font.Size = 13
$ws.Range("D1").Font.Name = "Helvetica Neue"

# ---------------------------------------------------------------------
# 3. Column widths
# ---------------------------------------------------------------------
$ws.Columns("B").ColumnWidth = 65.917
$ws.Columns("C").ColumnWidth = 59.584
$ws.Columns("C").BestFit = $true
$ws.Columns("D").ColumnWidth = 30.25

# ---------------------------------------------------------------------
# 4. Row 1 height
# ---------------------------------------------------------------------
$ws.Rows(1).RowHeight = 17

# ---------------------------------------------------------------------
# 5. Selection / active cell moves to B13
# ---------------------------------------------------------------------
[void]$ws.Range("B13").Select()

Write-Host "CapellaForm: added Locator Type column (D) with Xpath/CSS classification"
